$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New valve-state data rows (D/E/F/G/H 9:17) -- one-hot pattern that
#    cycles through the five valve outputs two rows at a time.
# ---------------------------------------------------------------------
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0

$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0

$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0

$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0

$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0

$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 1

$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1

# ---------------------------------------------------------------------
# 2. Extend the table with blank rows down to row 30 (values only,
#    formatting is applied to the whole block below).
# ---------------------------------------------------------------------
$ws.Range("C18:H30").ClearContents()

# ---------------------------------------------------------------------
# 3. Give the whole table (headers + data + blank rows) a thin box
#    border around every cell plus centred alignment.
# ---------------------------------------------------------------------
$table = $ws.Range("C7:H30")
$table.HorizontalAlignment = -4108
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2

# ---------------------------------------------------------------------
# 4. Colour every "closed" (=0) valve cell with a theme accent colour.
# ---------------------------------------------------------------------
$themeCells = "E8","F8","G8","H8","E9","F9","G9","H9","D10","F10","G10","H10","D11","F11","G11","H11","D12","E12","G12","H12","D13","E13","G13","H13","D14","E14","F14","H14","D15","E15","F15","H15","D16","E16","F16","G16","D17","E17","F17","G17"
foreach ($addr in $themeCells) {
    $ws.Range($addr).Interior.ThemeColor = 6
}

# ---------------------------------------------------------------------
# 5. Colour every "open" (=1) valve cell green.
# ---------------------------------------------------------------------
$greenCells = "D8","D9","E10","E11","F12","F13","G14","G15"
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.Color = 5296274
}

# ---------------------------------------------------------------------
# 6. The last two outputs (H16:H17) are newly wired -- highlight them
#    green too and force the explicit (non-theme) font.
# ---------------------------------------------------------------------
$ws.Range("H16:H17").Interior.Color = 5296274
$ws.Range("H16:H17").Font.Name = "Calibri"

# ---------------------------------------------------------------------
# 7. Misc. view / page bits.
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
[void]$ws.Range("A12").Select()

Write-Host "done"
